$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $origStyle = $range.Style
    $range.Value = "'" + $text
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "72.725.84"
$ws.Range("E2").Value = "  -0.31%  "

$ws.Range("D3").Value = "3.937.07"
$ws.Range("E3").Value = "  -2.22%  "

$ws.Range("E4").Value = "  -0.11%  "

Set-TextValue $ws.Range("D5") "602.40"
$ws.Range("E5").Value = "  +1.30%  "

Set-TextValue $ws.Range("D6") "172.40"
$ws.Range("E6").Value = "  +12.42%  "

$ws.Range("E7").Value = "  -0.56%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  +3.09%  "

Set-TextValue $ws.Range("D10") "0.186"
$ws.Range("E10").Value = "  +9.24%  "

Set-TextValue $ws.Range("D11") "55.99"
$ws.Range("E11").Value = "  +2.91%  "

Set-TextValue $ws.Range("D12") "0.0000332"
$ws.Range("E12").Value = "  +3.50%  "

Set-TextValue $ws.Range("D13") "11.53"
$ws.Range("E13").Value = "  +4.94%  "

$ws.Range("D14").Value = "4.554.35"
$ws.Range("E14").Value = "  -2.55%  "

Set-TextValue $ws.Range("D15") "21.51"
$ws.Range("E15").Value = "  +4.01%  "

$ws.Range("D16").Value = "3.915.80"
$ws.Range("E16").Value = "  -2.85%  "

Set-TextValue $ws.Range("D17") "14.11"
$ws.Range("E17").Value = "  -1.35%  "

$ws.Range("E18").Value = "  -3.04%  "

$ws.Range("D19").Value = "72.551.80"
$ws.Range("E19").Value = "  -0.40%  "

$ws.Range("E20").Value = "  -1.14%  "

$ws.Range("E21").Value = "  +0.75%  "

$ws.Range("E22").Value = "  -0.24%  "

Set-TextValue $ws.Range("D23") "95.55"
$ws.Range("E23").Value = "  -1.91%  "

Set-TextValue $ws.Range("D24") "3.33"
$ws.Range("E24").Value = "  -6.07%  "

Set-TextValue $ws.Range("D25") "14.05"
$ws.Range("E25").Value = "  -2.16%  "

Set-TextValue $ws.Range("D26") "4.31"
$ws.Range("E26").Value = "  -0.28%  "

Set-TextValue $ws.Range("D27") "11.12"
$ws.Range("E27").Value = "  -3.71%  "

Set-TextValue $ws.Range("D28") "5.94"
$ws.Range("E28").Value = "  +0.21%  "

Set-TextValue $ws.Range("D29") "10.40"
$ws.Range("E29").Value = "  -3.42%  "

Set-TextValue $ws.Range("D30") "35.83"
$ws.Range("E30").Value = "  -2.65%  "

Set-TextValue $ws.Range("D31") "7.85"
$ws.Range("E31").Value = "  -1.35%  "

Set-TextValue $ws.Range("D32") "13.83"
$ws.Range("E32").Value = "  +1.24%  "

Set-TextValue $ws.Range("D33") "50.55"
$ws.Range("E33").Value = "  +0.54%  "

$ws.Range("E34").Value = "  -4.31%  "

Set-TextValue $ws.Range("D35") "0.0000100"
$ws.Range("E35").Value = "  +14.38%  "

Set-TextValue $ws.Range("D36") "68.65"
$ws.Range("E36").Value = "  -4.23%  "

Set-TextValue $ws.Range("D37") "631.73"
$ws.Range("E37").Value = "  -8.59%  "

Set-TextValue $ws.Range("D38") "0.426"
$ws.Range("E38").Value = "  -4.80%  "

Set-TextValue $ws.Range("D39") "3.39"
$ws.Range("E39").Value = "  +0.69%  "

$ws.Range("E40").Value = "  -0.18%  "

Set-TextValue $ws.Range("D41") "0.145"
$ws.Range("E41").Value = "  -2.56%  "

$ws.Range("E42").Value = "  +0.20%  "

Set-TextValue $ws.Range("D43") "3.30"
$ws.Range("E43").Value = "  +43.69%  "

Set-TextValue $ws.Range("D44") "0.0477"
$ws.Range("E44").Value = "  -3.26%  "

Set-TextValue $ws.Range("D45") "10.54"
$ws.Range("E45").Value = "  -5.98%  "

$ws.Range("E46").Value = "  -2.54%  "

$ws.Range("E47").Value = "  -5.26%  "

$ws.Range("E48").Value = "  -0.38%  "

Set-TextValue $ws.Range("D49") "2.85"
$ws.Range("E49").Value = "  -16.94%  "

$ws.Range("E50").Value = "  +4.70%  "

$ws.Range("D51").Value = "2.824.38"
$ws.Range("E51").Value = "  +0.03%  "
